$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new "Expected (Hours)" column
$ws.Cells.Item(1, 9).Value = "Expected (Hours)"
$ws.Cells.Item(1, 9).Font.Bold = $true

# Match formatting/column width behavior of column F (closest achievable to 17.42578125)
$ws.Columns.Item(9).ColumnWidth = 16.6

# Fill in the Expected (Hours) formula for rows 2 through 49
$ws.Range("I2").Formula = "=((F2+(4*G2)+H2)/6)"
$ws.Range("I3:I49").Formula = "=((F3+(4*G3)+H3)/6)"

$wb.Application.Calculate()
